$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.683.83'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').Value = '3.257.39'
$ws.Range('E3').Value = '  +5.11%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.49%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.251.73'
$ws.Range('E8').Value = '  +5.10%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.41'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('E12').Value = '  +2.28%  '
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.59'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D15').Value = '3.795.07'
$ws.Range('E15').Value = '  +4.72%  '
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '3.252.09'
$ws.Range('E17').Value = '  +4.06%  '
$ws.Range('D18').Value = '63.733.82'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '478.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('E22').Value = '  +5.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.63%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.88%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.96'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.60%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  +1.24%  '
$ws.Range('B28').Value = 'FirstDigitalUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.59%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.85%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.11%  '
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.08'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('D38').Value = '0.0₃0722'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('E39').Value = '  +1.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '423.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.07%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.000.21'
$ws.Range('E41').Value = '  +5.24%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.42'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.17%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.77'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.111'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.23%  '
$ws.Range('E45').Value = '  +3.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.04%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.97'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.114'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.34%  '
